$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "python DS" sheet (sheet1): drop the hard-coded style on B2/B3, move the
#    selection from B4 to B3.
# ---------------------------------------------------------------------------
$wsDS = $wb.Worksheets.Item("python DS")
$wsDS.Range("B2:B3").ClearFormats()

# ---------------------------------------------------------------------------
# 2. "python PQ" sheet (sheet2): widen column A, select the whole sheet.
# ---------------------------------------------------------------------------
$wsPQ = $wb.Worksheets.Item("python PQ")
$wsPQ.Columns.Item(1).ColumnWidth = 226.66666666666666

# ---------------------------------------------------------------------------
# 3. "LoginValidData" sheet (sheet6): explicit column widths, move selection
#    from B2 to B3.
# ---------------------------------------------------------------------------
$wsLVD = $wb.Worksheets.Item("LoginValidData")
$wsLVD.Columns.Item(1).ColumnWidth = 15.333333333333334
$wsLVD.Columns.Item(2).ColumnWidth = 25.166666666666668

# ---------------------------------------------------------------------------
# 4. New sheet: SearchArray -- copied from "python PQ" (keeps the shared
#    strings / cell styles lined up with the source workbook) and then
#    trimmed down to the "search" example only.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPQ.Copy($null, $lastSheet)
$wsSearch = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSearch.Name = "SearchArray"

# stash the "square sorted" result-cell format (s=4) before we start deleting
# rows, so we can stamp it back onto the new blank row below.
$wsSearch.Range("B7").Copy()
$wsSearch.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSearch.Rows.Item(7).Delete()
$wsSearch.Rows.Item(6).Delete()
$wsSearch.Rows.Item(5).Delete()
$wsSearch.Rows.Item(3).Delete()
$wsSearch.Rows.Item(2).Delete()

$wsSearch.Range("A1:B1").ClearFormats()

$wsSearch.Range("A2").Copy()
$wsSearch.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSearch.Range("D1").Copy()
$wsSearch.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsSearch.Range("D1").Clear()

$wsSearch.Columns.Item(1).ColumnWidth = 226.66666666666666

# ---------------------------------------------------------------------------
# 5. New sheet: MaxConsecutive -- copied from "python PQ", trimmed down to
#    the "findMaxConsecutiveOnes" example only.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPQ.Copy($null, $lastSheet)
$wsMaxCons = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMaxCons.Name = "MaxConsecutive"

$wsMaxCons.Range("B7").Copy()
$wsMaxCons.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMaxCons.Rows.Item(7).Delete()
$wsMaxCons.Rows.Item(6).Delete()
$wsMaxCons.Rows.Item(4).Delete()
$wsMaxCons.Rows.Item(3).Delete()
$wsMaxCons.Rows.Item(2).Delete()

# Only rows 1-2 are left now, so row 10 is a safe scratch spot for the
# header-style (s=1) stash used by B3 below.
$wsMaxCons.Range("A1").Copy()
$wsMaxCons.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# row3: blank "code"-style cell (A3) + blank header-style cell (B3), custom height
$wsMaxCons.Range("A2").Copy()
$wsMaxCons.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsMaxCons.Range("D10").Copy()
$wsMaxCons.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsMaxCons.Range("D1").Copy()
$wsMaxCons.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMaxCons.Range("D1").Clear()
$wsMaxCons.Range("D10").Clear()

$wsMaxCons.Rows.Item(3).RowHeight = 17.6

$wsMaxCons.Columns.Item(1).ColumnWidth = 226.66666666666666

# ---------------------------------------------------------------------------
# 6. New sheet: FindEventNum -- copied from "python PQ", trimmed down to the
#    "findNumbers" example only.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPQ.Copy($null, $lastSheet)
$wsFindEven = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsFindEven.Name = "FindEventNum"

$wsFindEven.Rows.Item(7).Delete()
$wsFindEven.Rows.Item(5).Delete()
$wsFindEven.Rows.Item(4).Delete()
$wsFindEven.Rows.Item(3).Delete()
$wsFindEven.Rows.Item(2).Delete()

$wsFindEven.Columns.Item(1).ColumnWidth = 226.66666666666666

# ---------------------------------------------------------------------------
# 7. New sheet: SquareSorted -- copied from "python PQ", trimmed down to the
#    "sortedSquares" example only.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPQ.Copy($null, $lastSheet)
$wsSquareSorted = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSquareSorted.Name = "SquareSorted"

$wsSquareSorted.Rows.Item(6).Delete()
$wsSquareSorted.Rows.Item(5).Delete()
$wsSquareSorted.Rows.Item(4).Delete()
$wsSquareSorted.Rows.Item(3).Delete()
$wsSquareSorted.Rows.Item(2).Delete()

$wsSquareSorted.Columns.Item(1).ColumnWidth = 226.66666666666666

# ---------------------------------------------------------------------------
# 8. Selections / active sheet bookkeeping - done last so the view state
#    (tabSelected / activeTab / per-sheet selection rectangle) matches.
# ---------------------------------------------------------------------------
$wsDS.Activate()
$wsDS.Range("B3").Select()

$wsPQ.Activate()
$wsPQ.Cells.Select()

$wsLVD.Activate()
$wsLVD.Range("B3").Select()

$wsMaxCons.Activate()
$wsMaxCons.Range("A2:XFD3").Select()

$wsFindEven.Activate()
$wsFindEven.Range("A2:XFD3").Select()

$wsSquareSorted.Activate()
$wsSquareSorted.Range("A2:XFD3").Select()

$wsSearch.Activate()
$wsSearch.Range("A30").Select()
